$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its original text formatting (many values look numeric,
# e.g. "1.00", "25.17" - without this Excel would coerce them to numbers and drop
# trailing zeros / thousands-style dots).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '61.327.42'
$ws.Range("E2").Value = '  +0.14%  '

# Row 3
$ws.Range("D3").Value = '2.381.06'
$ws.Range("E3").Value = '  +0.12%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").Value = '549.70'
$ws.Range("E5").Value = '  +0.16%  '

# Row 6
$ws.Range("D6").Value = '139.22'
$ws.Range("E6").Value = '  -1.74%  '

# Row 7
$ws.Range("D7").Value = '0.999'

# Row 8
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  -1.27%  '

# Row 9
$ws.Range("D9").Value = '2.379.63'
$ws.Range("E9").Value = '  +0.09%  '

# Row 10
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  +1.89%  '

# Row 11
$ws.Range("E11").Value = '  +1.53%  '

# Row 12
$ws.Range("D12").Value = '5.34'
$ws.Range("E12").Value = '  +0.76%  '

# Row 13
$ws.Range("D13").Value = '0.349'
$ws.Range("E13").Value = '  +0.67%  '

# Row 14
$ws.Range("D14").Value = '25.17'
$ws.Range("E14").Value = '  -0.57%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000167'
$ws.Range("E15").Value = '  +1.29%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '61.226.21'
$ws.Range("E16").Value = '  +0.10%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.368.97'
$ws.Range("E17").Value = '  -0.52%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '10.95'
$ws.Range("E18").Value = '  +2.19%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '4.15'
$ws.Range("E19").Value = '  +0.84%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '321.03'
$ws.Range("E20").Value = '  +0.89%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.72'
$ws.Range("E21").Value = '  +0.78%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.04%  '

# Row 23
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '64.25'
$ws.Range("E23").Value = '  +0.72%  '

# Row 24
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = '1.72'
$ws.Range("E24").Value = '  -9.41%  '

# Row 25
$ws.Range("B25").Value = 'Aptos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D25").Value = '8.52'
$ws.Range("E25").Value = '  +4.06%  '

# Row 26
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '8.19'
$ws.Range("E26").Value = '  +1.44%  '

# Row 27
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '0.0₃0898'
$ws.Range("E27").Value = '  -2.80%  '

# Row 28
$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D28").Value = '509.15'
$ws.Range("E28").Value = '  -3.26%  '

# Row 29
$ws.Range("E29").Value = '  +2.94%  '

# Row 30
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '1.38'
$ws.Range("E30").Value = '  -3.08%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '1.83'
$ws.Range("E31").Value = '  -0.16%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '1.53'
$ws.Range("E32").Value = '  -2.56%  '

# Row 33
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.16%  '

# Row 34
$ws.Range("B34").Value = 'Stacks'
$ws.Range("C34").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D34").Value = '1.92'
$ws.Range("E34").Value = '  +3.26%  '

# Row 35
$ws.Range("D35").Value = '4.70'
$ws.Range("E35").Value = '  +0.80%  '

# Row 36
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '5.45'
$ws.Range("E36").Value = '  +0.12%  '

# Row 37
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = '0.379'
$ws.Range("E37").Value = '  +1.19%  '

# Row 38
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '18.53'
$ws.Range("E38").Value = '  +2.34%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '146.90'
$ws.Range("E39").Value = '  +5.04%  '

# Row 40
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.08%  '

# Row 41
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").Value = '41.21'
$ws.Range("E41").Value = '  +1.78%  '

# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '151.13'
$ws.Range("E42").Value = '  +7.67%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.17'
$ws.Range("E43").Value = '  +1.29%  '

# Row 44
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '3.60'
$ws.Range("E44").Value = '  -0.09%  '

# Row 45
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = '0.0522'
$ws.Range("E45").Value = '  +1.08%  '

# Row 46
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '19.41'
$ws.Range("E46").Value = '  -3.35%  '

# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.578'
$ws.Range("E47").Value = '  +0.58%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.0906'
$ws.Range("E48").Value = '  +0.14%  '

# Row 49
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0224'
$ws.Range("E49").Value = '  -0.22%  '

# Row 50
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '11.42'
$ws.Range("E50").Value = '  +0.38%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '16.75'
$ws.Range("E51").Value = '  +0.09%  '
